$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("FilesTab") holds the file-listing SQL query in column B.
$cell = $ws.Cells.Item(5, 2)
$old = $cell.Value2

# Remove the now-redundant "fd.file_source AS "File Source"" column from the
# SELECT list and re-indent the following FROM clause, matching the
# upstream edit to this query.
$search = "    smp.sample_id AS `"Sample ID`",`n    fd.file_source AS `"File Source`"`nFROM `n"
$replace = "    smp.sample_id AS `"Sample ID`"`n  FROM `n"

$new = $old.Replace($search, $replace)
$cell.Value2 = $new

# Writing the cell re-triggers wrap-text autofit on this long query; restore
# the row to its original (already-at-cap) height so only the text content
# differs from the source file.
$ws.Rows.Item(5).RowHeight = 409.5
